# Add the "HMI Internal" global-variable-generator sheet to the PLC
# global_variable_template workbook (see commit: "Add HMI_internal into
# global_variable_generator").
#
# The new sheet is appended after "Pump" and becomes the active sheet, a
# small variable table (base_addr / var_name / var_type / addr_offset) is
# filled in, and the helper column G gets a currency-style number format
# on the rows that hold "addressable" entries.

$wb = $excel.ActiveWorkbook

# --- create the sheet, positioned after the last existing sheet (Pump) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "HMI Internal"

# --- data rows first (2..7), header row (1) last -------------------------
# Writing in this order keeps the shared-string table in the same order as
# the source: PUMP_0_SELECTED, BIT, PUMP_1_SELECTED, PUMP_2_SELECTED,
# PUMP_CHECK, DUMMY_BIT, STRING, var_name, var_type.

$ws.Range("A2").Value = 1000
$ws.Range("B2").Value = "PUMP_0_SELECTED"
$ws.Range("C2").Value = "BIT"
$ws.Range("D2").Value = 1
$ws.Range("G2").NumberFormat = "$#,##0.00_);[Red]($#,##0.00)"

$ws.Range("B3").Value = "PUMP_1_SELECTED"
$ws.Range("C3").Value = "BIT"
$ws.Range("D3").Value = 1
$ws.Range("G3").NumberFormat = "$#,##0.00_);[Red]($#,##0.00)"

$ws.Range("B4").Value = "PUMP_2_SELECTED"
$ws.Range("C4").Value = "BIT"
$ws.Range("D4").Value = 1
$ws.Range("G4").NumberFormat = "$#,##0.00_);[Red]($#,##0.00)"

$ws.Range("B5").Value = "PUMP_CHECK"
$ws.Range("C5").Value = "WORD"
$ws.Range("D5").Value = 1
$ws.Range("G5").NumberFormat = "$#,##0_);[Red]($#,##0)"

$ws.Range("B6").Value = "DUMMY_BIT"
$ws.Range("C6").Value = "BIT"
$ws.Range("D6").Value = 1
$ws.Range("G6").NumberFormat = "$#,##0.00_);[Red]($#,##0.00)"

$ws.Range("B7").Value = "STRING"
$ws.Range("C7").Value = "WORD"
$ws.Range("D7").Value = 20
$ws.Range("G7").NumberFormat = "$#,##0_);[Red]($#,##0)"

$ws.Range("A1").Value = "base_addr"
$ws.Range("B1").Value = "var_name"
$ws.Range("C1").Value = "var_type"
$ws.Range("D1").Value = "addr_offset"

# --- column widths (best-fit look, matching the other generator sheets) --
$ws.Columns("A").ColumnWidth = 10.140625
$ws.Columns("B").ColumnWidth = 17.85546875
$ws.Columns("C").ColumnWidth = 8.7109375
$ws.Columns("D").ColumnWidth = 11.28515625

# --- leave the same selection / active-sheet state as the source file ----
$ws.Range("F5").Select() | Out-Null
